$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 84
$ws.Range("C3").Value = 64
$ws.Range("D3").Value = 98
$ws.Range("E3").Value = 100
$ws.Range("G3").Value = 92
$ws.Range("I3").Value = 146
$ws.Range("J3").Value = 149
$ws.Range("B6").Value = 274
$ws.Range("E6").Value = 313
$ws.Range("G6").Value = 354
$ws.Range("H6").Value = 320
$ws.Range("I6").Value = 381
$ws.Range("J6").Value = 293
$ws.Range("K6").Value = 370
$ws.Range("B7").Value = 375
$ws.Range("C7").Value = 463
$ws.Range("D7").Value = 475
$ws.Range("E7").Value = 471
$ws.Range("G7").Value = 514
$ws.Range("H7").Value = 506
$ws.Range("I7").Value = 632
$ws.Range("J7").Value = 544
$ws.Range("K7").Value = 650

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 29
$ws.Range("E7").Value = 38

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("G3").Value = 5
$ws.Range("I6").Value = 15
$ws.Range("G7").Value = 20
$ws.Range("I7").Value = 31

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("H2").Value = 1
$ws.Range("H6").Value = 4

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("B4").Value = 10
$ws.Range("K4").Value = 3
$ws.Range("B5").Value = 12
$ws.Range("K5").Value = 6

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 6
$ws.Range("H8").Value = 36
$ws.Range("C29").Value = 5
$ws.Range("E32").Value = 38
$ws.Range("B35").Value = 6
$ws.Range("G36").Value = 20
$ws.Range("I36").Value = 31
$ws.Range("H41").Value = 4
$ws.Range("J49").Value = 4
$ws.Range("G53").Value = 62
$ws.Range("J53").Value = 89
$ws.Range("E62").Value = 5
$ws.Range("E65").Value = 8
$ws.Range("I70").Value = 15
$ws.Range("G72").Value = 5
$ws.Range("I76").Value = 17
$ws.Range("H79").Value = 8
$ws.Range("B80").Value = 12
$ws.Range("K80").Value = 6
$ws.Range("D85").Value = 5
$ws.Range("B92").Value = 8
$ws.Range("J96").Value = 8
$ws.Range("B98").Value = 375
$ws.Range("C98").Value = 463
$ws.Range("D98").Value = 475
$ws.Range("E98").Value = 471
$ws.Range("G98").Value = 514
$ws.Range("H98").Value = 506
$ws.Range("I98").Value = 632
$ws.Range("J98").Value = 544
$ws.Range("K98").Value = 650

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J5").Value = 4
$ws.Range("J6").Value = 8

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("G6").Value = 43
$ws.Range("J6").Value = 46
$ws.Range("G7").Value = 62
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 6

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 8

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("H2").Value = 3
$ws.Range("H6").Value = 8

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 4
$ws.Range("I6").Value = 17

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("E3").Value = 2
$ws.Range("E6").Value = 8

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("E6").Value = 4
$ws.Range("E7").Value = 5

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("C3").Value = 2
$ws.Range("C6").Value = 5

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 5

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 6

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("D3").Value = 1
$ws.Range("D5").Value = 5

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I3").Value = 1
$ws.Range("I5").Value = 4

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I3").Value = 2
$ws.Range("I5").Value = 15

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("H5").Value = 29
$ws.Range("H6").Value = 36
